$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tigran UPI update (field_wbddh_dsttl_upi, row 5)
$ws.Range("B5").Value = 336592

# Jemi and Shandao leaving -> field_wbddh_collaborator_upi (row 11) now a single UPI number
$ws.Range("B11").Value = 511294

# Highlight the two changed rows in yellow, like a manual change-tracking mark
$ws.Range("A5").Interior.Color = 65535
$ws.Range("B5").Interior.Color = 65535
$ws.Range("D5").Interior.Color = 65535

$ws.Range("A11").Interior.Color = 65535
$ws.Range("B11").Interior.Color = 65535
$ws.Range("D11").Interior.Color = 65535

# Leave the same selection state recorded by Excel on save
$ws.Range("H25").Select()
